$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39: phone number A39 was stored as text "71277620"; convert it to a
# real number (matches the diff: t="inlineStr" -> t="n").
$ws.Range("A39").Value = 71277620

# New row 40: a new payment record for phone 71277620 (Cash, $100).
# Column A is kept as TEXT (leading apostrophe forces text, matching the
# t="inlineStr" cell in the diff), then style is reset back to "Normal" so
# we don't leave a stray quote-prefix format on the cell.
$ws.Range("A40").Value = "'71277620"
$ws.Range("A40").Style = "Normal"

# B40 / F40 are blank text cells in the diff (t="inlineStr" with no content)
# rather than missing cells, so write an empty text value (apostrophe makes
# it take the text path) and reset the style afterwards.
$ws.Range("B40").Value = "'"
$ws.Range("B40").Style = "Normal"

$ws.Range("C40").Value = "Cash"
$ws.Range("D40").Value = "2025-08-18T17:29:17"
$ws.Range("E40").Value = 100

$ws.Range("F40").Value = "'"
$ws.Range("F40").Style = "Normal"

$ws.Range("G40").Value = 100
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
